# Scheduled market-data refresh: update computed price/profit columns (H-N)
# on the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2000  # ALC!H43 (2788 -> 2000)
$ws.Cells.Item(43, 10).Value = 0  # ALC!J43 (3576 -> 0)
$ws.Cells.Item(43, 12).Value = 0  # ALC!L43 (3576 -> 0)
$ws.Cells.Item(43, 14).Value = ""  # ALC!N43 clear (was -3714)
$ws.Cells.Item(70, 8).Value = 0  # ALC!H70 (3500 -> 0)
$ws.Cells.Item(70, 10).Value = 0  # ALC!J70 (3500 -> 0)
$ws.Cells.Item(70, 12).Value = 0  # ALC!L70 (10500 -> 0)
$ws.Cells.Item(70, 14).Value = ""  # ALC!N70 clear (was -11040)
$ws.Cells.Item(73, 8).Value = 0  # ALC!H73 (3500 -> 0)
$ws.Cells.Item(73, 10).Value = 0  # ALC!J73 (3500 -> 0)
$ws.Cells.Item(73, 12).Value = 0  # ALC!L73 (10500 -> 0)
$ws.Cells.Item(73, 14).Value = ""  # ALC!N73 clear (was -12372)
$ws.Cells.Item(98, 8).Value = 28662  # ALC!H98 (31197 -> 28662)
$ws.Cells.Item(98, 10).Value = 0  # ALC!J98 (34999.5 -> 0)
$ws.Cells.Item(98, 12).Value = 0  # ALC!L98 (34999.5 -> 0)
$ws.Cells.Item(98, 14).Value = ""  # ALC!N98 clear (was -37995.5)
$ws.Cells.Item(122, 8).Value = 28662  # ALC!H122 (31197 -> 28662)
$ws.Cells.Item(122, 10).Value = 0  # ALC!J122 (34999.5 -> 0)
$ws.Cells.Item(122, 12).Value = 0  # ALC!L122 (104998.5 -> 0)
$ws.Cells.Item(122, 14).Value = ""  # ALC!N122 clear (was -109898.5)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 5071  # ARM!H122 (5428.143 -> 5071)
$ws.Cells.Item(122, 9).Value = 4374.25  # ARM!I122 (4999.25 -> 4374.25)
$ws.Cells.Item(122, 11).Value = 13122.75  # ARM!K122 (14997.75 -> 13122.75)
$ws.Cells.Item(122, 13).Value = -10672.75  # ARM!M122 (-12547.75 -> -10672.75)
$ws.Cells.Item(140, 8).Value = 0  # ARM!H140 (50000 -> 0)
$ws.Cells.Item(140, 10).Value = 0  # ARM!J140 (50000 -> 0)
$ws.Cells.Item(140, 12).Value = 0  # ARM!L140 (50000 -> 0)
$ws.Cells.Item(140, 14).Value = ""  # ARM!N140 clear (was -60360)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 749.5  # BSM!H22 (749.8 -> 749.5)
$ws.Cells.Item(22, 9).Value = 749.5  # BSM!I22 (749.8 -> 749.5)
$ws.Cells.Item(22, 11).Value = 749.5  # BSM!K22 (749.8 -> 749.5)
$ws.Cells.Item(22, 13).Value = -576.5  # BSM!M22 (-576.8 -> -576.5)
$ws.Cells.Item(94, 8).Value = 2282.1177  # BSM!H94 (2374.8125 -> 2282.1177)
$ws.Cells.Item(94, 9).Value = 1779.7  # BSM!I94 (1888.6666 -> 1779.7)
$ws.Cells.Item(94, 11).Value = 1779.7  # BSM!K94 (1888.6666 -> 1779.7)
$ws.Cells.Item(94, 13).Value = -1328.7  # BSM!M94 (-1437.6666 -> -1328.7)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 24.928572  # CRP!H7 (25.428572 -> 24.928572)
$ws.Cells.Item(7, 9).Value = 22  # CRP!I7 (22.636364 -> 22)
$ws.Cells.Item(7, 11).Value = 22  # CRP!K7 (22.636364 -> 22)
$ws.Cells.Item(7, 13).Value = 91  # CRP!M7 (90.363636 -> 91)
$ws.Cells.Item(22, 8).Value = 754.4286  # CRP!H22 (686.55554 -> 754.4286)
$ws.Cells.Item(22, 9).Value = 754.4286  # CRP!I22 (686.55554 -> 754.4286)
$ws.Cells.Item(22, 11).Value = 754.4286  # CRP!K22 (686.55554 -> 754.4286)
$ws.Cells.Item(22, 13).Value = -404.4286  # CRP!M22 (-336.55554 -> -404.4286)
$ws.Cells.Item(31, 8).Value = 36579  # CRP!H31 (39197.273 -> 36579)
$ws.Cells.Item(31, 9).Value = 23732  # CRP!I31 (28238.6 -> 23732)
$ws.Cells.Item(31, 10).Value = 51995.4  # CRP!J31 (48329.5 -> 51995.4)
$ws.Cells.Item(31, 11).Value = 23732  # CRP!K31 (28238.6 -> 23732)
$ws.Cells.Item(31, 12).Value = 51995.4  # CRP!L31 (48329.5 -> 51995.4)
$ws.Cells.Item(31, 13).Value = -23437  # CRP!M31 (-27943.6 -> -23437)
$ws.Cells.Item(31, 14).Value = -52585.4  # CRP!N31 (-48919.5 -> -52585.4)
$ws.Cells.Item(34, 8).Value = 36579  # CRP!H34 (39197.273 -> 36579)
$ws.Cells.Item(34, 9).Value = 23732  # CRP!I34 (28238.6 -> 23732)
$ws.Cells.Item(34, 10).Value = 51995.4  # CRP!J34 (48329.5 -> 51995.4)
$ws.Cells.Item(34, 11).Value = 23732  # CRP!K34 (28238.6 -> 23732)
$ws.Cells.Item(34, 12).Value = 51995.4  # CRP!L34 (48329.5 -> 51995.4)
$ws.Cells.Item(34, 13).Value = -23530  # CRP!M34 (-28036.6 -> -23530)
$ws.Cells.Item(34, 14).Value = -52399.4  # CRP!N34 (-48733.5 -> -52399.4)
$ws.Cells.Item(86, 8).Value = 9290.666999999999  # CRP!H86 (8849.143 -> 9290.666999999999)
$ws.Cells.Item(86, 9).Value = 9759  # CRP!I86 (9165.833000000001 -> 9759)
$ws.Cells.Item(86, 11).Value = 9759  # CRP!K86 (9165.833000000001 -> 9759)
$ws.Cells.Item(86, 13).Value = -8636  # CRP!M86 (-8042.833000000001 -> -8636)
$ws.Cells.Item(89, 8).Value = 9290.666999999999  # CRP!H89 (8849.143 -> 9290.666999999999)
$ws.Cells.Item(89, 9).Value = 9759  # CRP!I89 (9165.833000000001 -> 9759)
$ws.Cells.Item(89, 11).Value = 48795  # CRP!K89 (45829.165 -> 48795)
$ws.Cells.Item(89, 13).Value = -43179  # CRP!M89 (-40213.165 -> -43179)
$ws.Cells.Item(94, 8).Value = 1859.5  # CRP!H94 (1708 -> 1859.5)
$ws.Cells.Item(94, 9).Value = 865.6667  # CRP!I94 (849 -> 865.6667)
$ws.Cells.Item(94, 11).Value = 865.6667  # CRP!K94 (849 -> 865.6667)
$ws.Cells.Item(94, 13).Value = -414.6667  # CRP!M94 (-398 -> -414.6667)
$ws.Cells.Item(99, 8).Value = 2001187.5  # CRP!H99 (1144507 -> 2001187.5)
$ws.Cells.Item(99, 9).Value = 2501000  # CRP!I99 (1251499.8 -> 2501000)
$ws.Cells.Item(99, 10).Value = 1501375  # CRP!J99 (1001850 -> 1501375)
$ws.Cells.Item(99, 11).Value = 2501000  # CRP!K99 (1251499.8 -> 2501000)
$ws.Cells.Item(99, 12).Value = 1501375  # CRP!L99 (1001850 -> 1501375)
$ws.Cells.Item(99, 13).Value = -2499502  # CRP!M99 (-1250001.8 -> -2499502)
$ws.Cells.Item(99, 14).Value = -1504371  # CRP!N99 (-1004846 -> -1504371)
$ws.Cells.Item(107, 8).Value = 190.05882  # CRP!H107 (198.875 -> 190.05882)
$ws.Cells.Item(107, 9).Value = 210.3077  # CRP!I107 (223.75 -> 210.3077)
$ws.Cells.Item(107, 11).Value = 210.3077  # CRP!K107 (223.75 -> 210.3077)
$ws.Cells.Item(107, 13).Value = 1709.6923  # CRP!M107 (1696.25 -> 1709.6923)
$ws.Cells.Item(122, 8).Value = 1507.125  # CRP!H122 (1221.0714 -> 1507.125)
$ws.Cells.Item(122, 9).Value = 1508.1428  # CRP!I122 (1284.5 -> 1508.1428)
$ws.Cells.Item(122, 10).Value = 1500  # CRP!J122 (1062.5 -> 1500)
$ws.Cells.Item(122, 11).Value = 4524.428400000001  # CRP!K122 (3853.5 -> 4524.428400000001)
$ws.Cells.Item(122, 12).Value = 4500  # CRP!L122 (3187.5 -> 4500)
$ws.Cells.Item(122, 13).Value = -2074.428400000001  # CRP!M122 (-1403.5 -> -2074.428400000001)
$ws.Cells.Item(122, 14).Value = -9400  # CRP!N122 (-8087.5 -> -9400)
$ws.Cells.Item(126, 8).Value = 2001187.5  # CRP!H126 (1144507 -> 2001187.5)
$ws.Cells.Item(126, 9).Value = 2501000  # CRP!I126 (1251499.8 -> 2501000)
$ws.Cells.Item(126, 10).Value = 1501375  # CRP!J126 (1001850 -> 1501375)
$ws.Cells.Item(126, 11).Value = 7503000  # CRP!K126 (3754499.4 -> 7503000)
$ws.Cells.Item(126, 12).Value = 4504125  # CRP!L126 (3005550 -> 4504125)
$ws.Cells.Item(126, 13).Value = -7500530  # CRP!M126 (-3752029.4 -> -7500530)
$ws.Cells.Item(126, 14).Value = -4509065  # CRP!N126 (-3010490 -> -4509065)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 4683.8  # CUL!H14 (4236.5 -> 4683.8)
$ws.Cells.Item(14, 9).Value = 4683.8  # CUL!I14 (4236.5 -> 4683.8)
$ws.Cells.Item(14, 11).Value = 14051.4  # CUL!K14 (12709.5 -> 14051.4)
$ws.Cells.Item(14, 13).Value = -13878.4  # CUL!M14 (-12536.5 -> -13878.4)
$ws.Cells.Item(86, 8).Value = 1791.6  # CUL!H86 (1678 -> 1791.6)
$ws.Cells.Item(86, 10).Value = 1791.6  # CUL!J86 (1678 -> 1791.6)
$ws.Cells.Item(86, 12).Value = 5374.799999999999  # CUL!L86 (5034 -> 5374.799999999999)
$ws.Cells.Item(86, 14).Value = -7746.799999999999  # CUL!N86 (-7406 -> -7746.799999999999)
$ws.Cells.Item(89, 8).Value = 1791.6  # CUL!H89 (1678 -> 1791.6)
$ws.Cells.Item(89, 10).Value = 1791.6  # CUL!J89 (1678 -> 1791.6)
$ws.Cells.Item(89, 12).Value = 16124.4  # CUL!L89 (15102 -> 16124.4)
$ws.Cells.Item(89, 14).Value = -27980.4  # CUL!N89 (-26958 -> -27980.4)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7915.3335  # LTW!H7 (8248.25 -> 7915.3335)
$ws.Cells.Item(7, 9).Value = 7998.6  # LTW!I7 (8248.25 -> 7998.6)
$ws.Cells.Item(7, 10).Value = 7499  # LTW!J7 (0 -> 7499)
$ws.Cells.Item(7, 11).Value = 7998.6  # LTW!K7 (8248.25 -> 7998.6)
$ws.Cells.Item(7, 12).Value = 7499  # LTW!L7 (0 -> 7499)
$ws.Cells.Item(7, 13).Value = -7886.6  # LTW!M7 (-8136.25 -> -7886.6)
$ws.Cells.Item(7, 14).Value = -7723  # LTW!N7 (None -> -7723)
$ws.Cells.Item(61, 8).Value = 4509.5454  # LTW!H61 (5240.5 -> 4509.5454)
$ws.Cells.Item(61, 9).Value = 4825  # LTW!I61 (5205.5557 -> 4825)
$ws.Cells.Item(61, 10).Value = 1355  # LTW!J61 (5555 -> 1355)
$ws.Cells.Item(61, 11).Value = 4825  # LTW!K61 (5205.5557 -> 4825)
$ws.Cells.Item(61, 12).Value = 1355  # LTW!L61 (5555 -> 1355)
$ws.Cells.Item(61, 13).Value = -4623  # LTW!M61 (-5003.5557 -> -4623)
$ws.Cells.Item(61, 14).Value = -1759  # LTW!N61 (-5959 -> -1759)
$ws.Cells.Item(68, 10).Value = 2000  # LTW!J68 (0 -> 2000)
$ws.Cells.Item(68, 12).Value = 2000  # LTW!L68 (0 -> 2000)
$ws.Cells.Item(68, 14).Value = -3498  # LTW!N68 (None -> -3498)
$ws.Cells.Item(71, 10).Value = 2000  # LTW!J71 (0 -> 2000)
$ws.Cells.Item(71, 12).Value = 10000  # LTW!L71 (0 -> 10000)
$ws.Cells.Item(71, 14).Value = -17488  # LTW!N71 (None -> -17488)
$ws.Cells.Item(82, 8).Value = 1600  # LTW!H82 (1611.1111 -> 1600)
$ws.Cells.Item(82, 10).Value = 1833.3334  # LTW!J82 (2000 -> 1833.3334)
$ws.Cells.Item(82, 12).Value = 1833.3334  # LTW!L82 (2000 -> 1833.3334)
$ws.Cells.Item(82, 14).Value = -2555.3334  # LTW!N82 (-2722 -> -2555.3334)
$ws.Cells.Item(85, 8).Value = 1600  # LTW!H85 (1611.1111 -> 1600)
$ws.Cells.Item(85, 10).Value = 1833.3334  # LTW!J85 (2000 -> 1833.3334)
$ws.Cells.Item(85, 12).Value = 1833.3334  # LTW!L85 (2000 -> 1833.3334)
$ws.Cells.Item(85, 14).Value = -4329.3334  # LTW!N85 (-4496 -> -4329.3334)
$ws.Cells.Item(113, 8).Value = 4509.5454  # LTW!H113 (5240.5 -> 4509.5454)
$ws.Cells.Item(113, 9).Value = 4825  # LTW!I113 (5205.5557 -> 4825)
$ws.Cells.Item(113, 10).Value = 1355  # LTW!J113 (5555 -> 1355)
$ws.Cells.Item(113, 11).Value = 4825  # LTW!K113 (5205.5557 -> 4825)
$ws.Cells.Item(113, 12).Value = 1355  # LTW!L113 (5555 -> 1355)
$ws.Cells.Item(113, 13).Value = -2655  # LTW!M113 (-3035.5557 -> -2655)
$ws.Cells.Item(113, 14).Value = -5695  # LTW!N113 (-9895 -> -5695)
$ws.Cells.Item(122, 8).Value = 2499  # LTW!H122 (2750 -> 2499)
$ws.Cells.Item(122, 9).Value = 1998.75  # LTW!I122 (2000 -> 1998.75)
$ws.Cells.Item(122, 10).Value = 3499.5  # LTW!J122 (5000 -> 3499.5)
$ws.Cells.Item(122, 11).Value = 5996.25  # LTW!K122 (6000 -> 5996.25)
$ws.Cells.Item(122, 12).Value = 10498.5  # LTW!L122 (15000 -> 10498.5)
$ws.Cells.Item(122, 13).Value = -3546.25  # LTW!M122 (-3550 -> -3546.25)
$ws.Cells.Item(122, 14).Value = -15398.5  # LTW!N122 (-19900 -> -15398.5)
$ws.Cells.Item(126, 8).Value = 7915.3335  # LTW!H126 (8248.25 -> 7915.3335)
$ws.Cells.Item(126, 9).Value = 7998.6  # LTW!I126 (8248.25 -> 7998.6)
$ws.Cells.Item(126, 10).Value = 7499  # LTW!J126 (0 -> 7499)
$ws.Cells.Item(126, 11).Value = 23995.8  # LTW!K126 (24744.75 -> 23995.8)
$ws.Cells.Item(126, 12).Value = 22497  # LTW!L126 (0 -> 22497)
$ws.Cells.Item(126, 13).Value = -21525.8  # LTW!M126 (-22274.75 -> -21525.8)
$ws.Cells.Item(126, 14).Value = -27437  # LTW!N126 (None -> -27437)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 1917.1666  # WVR!H4 (903 -> 1917.1666)
$ws.Cells.Item(4, 9).Value = 3500  # WVR!I4 (0 -> 3500)
$ws.Cells.Item(4, 10).Value = 1600.6  # WVR!J4 (903 -> 1600.6)
$ws.Cells.Item(4, 11).Value = 3500  # WVR!K4 (0 -> 3500)
$ws.Cells.Item(4, 12).Value = 1600.6  # WVR!L4 (903 -> 1600.6)
$ws.Cells.Item(4, 13).Value = -3387  # WVR!M4 (None -> -3387)
$ws.Cells.Item(4, 14).Value = -1826.6  # WVR!N4 (-1129 -> -1826.6)
$ws.Cells.Item(62, 8).Value = 0  # WVR!H62 (2500 -> 0)
$ws.Cells.Item(62, 10).Value = 0  # WVR!J62 (2500 -> 0)
$ws.Cells.Item(62, 12).Value = 0  # WVR!L62 (2500 -> 0)
$ws.Cells.Item(62, 14).Value = ""  # WVR!N62 clear (was -3748)
$ws.Cells.Item(65, 8).Value = 0  # WVR!H65 (2500 -> 0)
$ws.Cells.Item(65, 10).Value = 0  # WVR!J65 (2500 -> 0)
$ws.Cells.Item(65, 12).Value = 0  # WVR!L65 (12500 -> 0)
$ws.Cells.Item(65, 14).Value = ""  # WVR!N65 clear (was -18740)
$ws.Cells.Item(122, 8).Value = 1012.25  # WVR!H122 (987.7778 -> 1012.25)
$ws.Cells.Item(122, 9).Value = 919.8  # WVR!I122 (898.5 -> 919.8)
$ws.Cells.Item(122, 11).Value = 2759.4  # WVR!K122 (2695.5 -> 2759.4)
$ws.Cells.Item(122, 13).Value = -309.3999999999996  # WVR!M122 (-245.5 -> -309.3999999999996)
$ws.Cells.Item(126, 8).Value = 2902.875  # WVR!H126 (3027.25 -> 2902.875)
$ws.Cells.Item(126, 9).Value = 2398.8  # WVR!I126 (2597.8 -> 2398.8)
$ws.Cells.Item(126, 11).Value = 7196.400000000001  # WVR!K126 (7793.400000000001 -> 7196.400000000001)
$ws.Cells.Item(126, 13).Value = -4726.400000000001  # WVR!M126 (-5323.400000000001 -> -4726.400000000001)
